# Update cryptocurrency price/volume snapshot values (Price = column D,
# Volume(1h) = column E) to match the latest scrape.
#
# Values are written with a leading apostrophe so Excel keeps them as
# literal text (matching the original inlineStr/text cells) instead of
# re-interpreting numeric- or percent-looking strings as actual numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.31"
$ws.Range("E2").Value = "'-2.81%"

$ws.Range("D3").Value = "'42.69"
$ws.Range("E3").Value = "'-5.91%"

$ws.Range("D4").Value = "'5.174"
$ws.Range("E4").Value = "'-7.32%"

$ws.Range("D5").Value = "'0.08184"
$ws.Range("E5").Value = "'-1.90%"

$ws.Range("D6").Value = "'4.301"
$ws.Range("E6").Value = "'-3.32%"

$ws.Range("D7").Value = "'1.808"
$ws.Range("E7").Value = "'-14.11%"

$ws.Range("D8").Value = "'0.9340"
$ws.Range("E8").Value = "'-4.64%"

$ws.Range("D9").Value = "'0.1109"
$ws.Range("E9").Value = "'-7.35%"

$ws.Range("D10").Value = "'0.1863"
$ws.Range("E10").Value = "'-3.05%"

$ws.Range("D11").Value = "'0.09432"
$ws.Range("E11").Value = "'-4.75%"

$ws.Range("D12").Value = "'0.04687"
$ws.Range("E12").Value = "'1.31%"

$ws.Range("D13").Value = "'7.418"
$ws.Range("E13").Value = "'-28.08%"

$ws.Range("D14").Value = "'0.1058"
$ws.Range("E14").Value = "'0.07%"

$ws.Range("D15").Value = "'0.001305"
$ws.Range("E15").Value = "'2.13%"

$ws.Range("D16").Value = "'0.005893"
$ws.Range("E16").Value = "'-0.35%"

$ws.Range("D17").Value = "'3.353"
$ws.Range("E17").Value = "'-0.67%"

$ws.Range("D18").Value = "'2.504"
$ws.Range("E18").Value = "'-2.83%"

$ws.Range("D19").Value = "'0.3379"
$ws.Range("E19").Value = "'1.11%"

$ws.Range("E20").Value = "'-0.09%"

$ws.Range("D21").Value = "'0.2548"
$ws.Range("E21").Value = "'-8.47%"

$ws.Range("D22").Value = "'0.04148"
$ws.Range("E22").Value = "'-0.79%"

$ws.Range("D23").Value = "'0.001250"
$ws.Range("E23").Value = "'-3.40%"

$ws.Range("E24").Value = "'-5.04%"

$ws.Range("D25").Value = "'0.0001202"
$ws.Range("E25").Value = "'-7.71%"

$ws.Range("D26").Value = "'0.0002983"
$ws.Range("E26").Value = "'-20.38%"

$ws.Range("D38").Value = "'0.02738"
$ws.Range("E38").Value = "'1.30%"

$ws.Range("D39").Value = "'0.05575"
$ws.Range("E39").Value = "'-2.96%"

$ws.Range("D40").Value = "'0.008036"
$ws.Range("E40").Value = "'1.85%"

$ws.Range("D41").Value = "'0.1399"
$ws.Range("E41").Value = "'-2.24%"

$ws.Range("D42").Value = "'0.006558"
$ws.Range("E42").Value = "'-12.79%"

$ws.Range("D43").Value = "'0.002068"
$ws.Range("E43").Value = "'-1.50%"

$ws.Range("D44").Value = "'0.008281"
$ws.Range("E44").Value = "'-2.74%"

$ws.Range("D45").Value = "'0.3491"
$ws.Range("E45").Value = "'3.55%"

$ws.Range("D46").Value = "'0.00006944"
$ws.Range("E46").Value = "'-2.55%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.05%"

$ws.Range("D48").Value = "'0.003513"
$ws.Range("E48").Value = "'-0.41%"

$ws.Range("D49").Value = "'0.003535"
$ws.Range("E49").Value = "'-0.04%"

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'-0.05%"

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.05%"
